$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1387.4
$ws.Cells.Item(28, 9).Value = 557.2353000000001
$ws.Cells.Item(28, 11).Value = 557.2353000000001
$ws.Cells.Item(28, 13).Value = -72.23530000000005

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 2069.5
$ws.Cells.Item(112, 10).Value = 2069.5
$ws.Cells.Item(112, 12).Value = 6208.5
$ws.Cells.Item(112, 14).Value = -8424.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 4514.5
$ws.Cells.Item(135, 9).Value = 3995.5
$ws.Cells.Item(135, 11).Value = 35959.5
$ws.Cells.Item(135, 13).Value = -33424.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 3314.8772
$ws.Cells.Item(138, 9).Value = 4306.353
$ws.Cells.Item(138, 10).Value = 2893.5
$ws.Cells.Item(138, 11).Value = 12919.059
$ws.Cells.Item(138, 12).Value = 8680.5
$ws.Cells.Item(138, 13).Value = -7779.059000000001
$ws.Cells.Item(138, 14).Value = -18960.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 735
$ws.Cells.Item(2, 9).Value = 455.05884
$ws.Cells.Item(2, 11).Value = 455.05884
$ws.Cells.Item(2, 13).Value = -342.05884

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3086.4614
$ws.Cells.Item(61, 9).Value = 2249.1052
$ws.Cells.Item(61, 11).Value = 2249.1052
$ws.Cells.Item(61, 13).Value = -2037.1052

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 4257.643
$ws.Cells.Item(63, 9).Value = 4488.1665
$ws.Cells.Item(63, 11).Value = 4488.1665
$ws.Cells.Item(63, 13).Value = -3802.1665

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 4257.643
$ws.Cells.Item(66, 9).Value = 4488.1665
$ws.Cells.Item(66, 11).Value = 22440.8325
$ws.Cells.Item(66, 13).Value = -19008.8325

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(92, 8).Value = 24998
$ws.Cells.Item(92, 10).Value = 24998
$ws.Cells.Item(92, 12).Value = 24998
$ws.Cells.Item(92, 14).Value = -29990

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 735
$ws.Cells.Item(116, 9).Value = 455.05884
$ws.Cells.Item(116, 11).Value = 455.05884
$ws.Cells.Item(116, 13).Value = 1838.94116

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1829.7441
$ws.Cells.Item(132, 9).Value = 1022.25
$ws.Cells.Item(132, 10).Value = 4178.8184
$ws.Cells.Item(132, 11).Value = 3066.75
$ws.Cells.Item(132, 12).Value = 12536.4552
$ws.Cells.Item(132, 13).Value = -536.75
$ws.Cells.Item(132, 14).Value = -17596.4552

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 3086.4614
$ws.Cells.Item(136, 9).Value = 2249.1052
$ws.Cells.Item(136, 11).Value = 6747.3156
$ws.Cells.Item(136, 13).Value = -4197.3156

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 735
$ws.Cells.Item(3, 9).Value = 455.05884
$ws.Cells.Item(3, 11).Value = 455.05884
$ws.Cells.Item(3, 13).Value = -341.05884

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 7900.5386
$ws.Cells.Item(99, 9).Value = 3873.25
$ws.Cells.Item(99, 11).Value = 3873.25
$ws.Cells.Item(99, 13).Value = -2375.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(100, 8).Value = 50999
$ws.Cells.Item(100, 10).Value = 50999
$ws.Cells.Item(100, 12).Value = 50999
$ws.Cells.Item(100, 14).Value = -53163

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3249.566
$ws.Cells.Item(31, 9).Value = 2041
$ws.Cells.Item(31, 10).Value = 8446.4
$ws.Cells.Item(31, 11).Value = 2041
$ws.Cells.Item(31, 12).Value = 8446.4
$ws.Cells.Item(31, 13).Value = -1746
$ws.Cells.Item(31, 14).Value = -9036.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3249.566
$ws.Cells.Item(34, 9).Value = 2041
$ws.Cells.Item(34, 10).Value = 8446.4
$ws.Cells.Item(34, 11).Value = 2041
$ws.Cells.Item(34, 12).Value = 8446.4
$ws.Cells.Item(34, 13).Value = -1839
$ws.Cells.Item(34, 14).Value = -8850.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 1856.8182
$ws.Cells.Item(33, 10).Value = 2037.5
$ws.Cells.Item(33, 12).Value = 12225
$ws.Cells.Item(33, 14).Value = -12791

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(60, 8).Value = 1114744.6
$ws.Cells.Item(60, 9).Value = 3334837.2
$ws.Cells.Item(60, 11).Value = 10004511.6
$ws.Cells.Item(60, 13).Value = -10004260.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 659.6667
$ws.Cells.Item(107, 10).Value = 659.6667
$ws.Cells.Item(107, 12).Value = 1979.0001
$ws.Cells.Item(107, 14).Value = -5819.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 428918.44
$ws.Cells.Item(121, 10).Value = 333933.34
$ws.Cells.Item(121, 12).Value = 1001800.02
$ws.Cells.Item(121, 14).Value = -1004420.02

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2999.75
$ws.Cells.Item(102, 9).Value = 2999.5
$ws.Cells.Item(102, 10).Value = 3000
$ws.Cells.Item(102, 11).Value = 2999.5
$ws.Cells.Item(102, 12).Value = 3000
$ws.Cells.Item(102, 13).Value = -1377.5
$ws.Cells.Item(102, 14).Value = -6244

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 65999.336
$ws.Cells.Item(123, 10).Value = 89998.664
$ws.Cells.Item(123, 12).Value = 89998.664
$ws.Cells.Item(123, 14).Value = -94898.664

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3995.4614
$ws.Cells.Item(126, 9).Value = 2794.3
$ws.Cells.Item(126, 10).Value = 7999.3335
$ws.Cells.Item(126, 11).Value = 8382.900000000001
$ws.Cells.Item(126, 12).Value = 23998.0005
$ws.Cells.Item(126, 13).Value = -5912.900000000001
$ws.Cells.Item(126, 14).Value = -28938.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3733.8823
$ws.Cells.Item(132, 9).Value = 3325.7
$ws.Cells.Item(132, 10).Value = 4317
$ws.Cells.Item(132, 11).Value = 9977.099999999999
$ws.Cells.Item(132, 12).Value = 12951
$ws.Cells.Item(132, 13).Value = -7447.099999999999
$ws.Cells.Item(132, 14).Value = -18011

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 13).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 17606.857
$ws.Cells.Item(41, 9).Value = 13386
$ws.Cells.Item(41, 10).Value = 19295.2
$ws.Cells.Item(41, 11).Value = 13386
$ws.Cells.Item(41, 12).Value = 19295.2
$ws.Cells.Item(41, 13).Value = -12996
$ws.Cells.Item(41, 14).Value = -20075.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 44544.4
$ws.Cells.Item(45, 10).Value = 44544.4
$ws.Cells.Item(45, 12).Value = 44544.4
$ws.Cells.Item(45, 14).Value = -45526.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3491.3235
$ws.Cells.Item(132, 10).Value = 2690.3333
$ws.Cells.Item(132, 12).Value = 8070.999899999999
$ws.Cells.Item(132, 14).Value = -13130.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(135, 8).Value = 51238
$ws.Cells.Item(135, 10).Value = 51238
$ws.Cells.Item(135, 12).Value = 51238
$ws.Cells.Item(135, 14).Value = -61378

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 8307.25
$ws.Cells.Item(136, 9).Value = 2350.4285
$ws.Cells.Item(136, 10).Value = 50005
$ws.Cells.Item(136, 11).Value = 7051.2855
$ws.Cells.Item(136, 12).Value = 150015
$ws.Cells.Item(136, 13).Value = -4501.2855
$ws.Cells.Item(136, 14).Value = -155115

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(137, 8).Value = 79999.11
$ws.Cells.Item(137, 10).Value = 79999.11
$ws.Cells.Item(137, 12).Value = 79999.11
$ws.Cells.Item(137, 14).Value = -90199.11

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(139, 8).Value = 81249.125
$ws.Cells.Item(139, 10).Value = 79999
$ws.Cells.Item(139, 12).Value = 79999
$ws.Cells.Item(139, 14).Value = -90279

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(140, 8).Value = 148974.5
$ws.Cells.Item(140, 10).Value = 148974.5
$ws.Cells.Item(140, 12).Value = 148974.5
$ws.Cells.Item(140, 14).Value = -159334.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(141, 8).Value = 68788.52
$ws.Cells.Item(141, 10).Value = 68788.52
$ws.Cells.Item(141, 12).Value = 68788.52
$ws.Cells.Item(141, 14).Value = -79148.52
